# Linear Regression.xlsx - "Values" sheet tidy-up.
#
# The "x" label used to be stored with a trailing space ("x ") and the
# shared-string table had it first; re-enter the four header/label cells
# with their correct (trimmed) text so the workbook is re-saved with
# "x" (no trailing space) as its own shared string and the cells keep
# showing the same visible text as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "x"
$ws.Range("B1").Value = "y"
$ws.Range("A2").Value = "(independent variable)"
$ws.Range("B2").Value = "(dependent variable)"

# Reset the cursor back to the top-left of the scrollable (frozen) pane
# instead of leaving it parked on the stale "B6" reference.
[void]$ws.Range("A1").Select()
